$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the "Conversión del día" rates text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.15 = 11829.97 pesos`n✅ 11829.97 pesos = 3.14 = 977.71 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- tasas!N10, O10, N12, O12: updated rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 317
$ws2.Range("O10").Value = 3750.1
$ws2.Range("N12").Value = 3763
$ws2.Range("O12").Value = 311
